$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 4 with the second test mail entry
$ws.Range("A4").Value = "Is er al nieuws?"
$ws.Range("B4").Value = "mailmind.test@zohomail.eu"
$ws.Range("C4").Value = "Testmail #2: Is er al nieuws?"
$ws.Range("D4").Value = "Opvolging / Status"
$ws.Range("E4").Value = "Dank voor je bericht. We hebben je eerdere e-mail ontvangen en doorgestuurd naar klantenservice@bedrijf.nl."
$ws.Range("F4").Value = "2025-08-04 20:01:22"
$ws.Range("G4").Value = "Ja"
$ws.Range("H4").Value = "Ja"
$ws.Range("I4").Value = "Nee"
$ws.Range("J4").Value = "Nee"

# Extend conditional formatting ranges to include the new row
foreach ($fc in $ws.Range("D2:D3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("D2:D4"))
}
foreach ($fc in $ws.Range("G2:G3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("G2:G4"))
}
foreach ($fc in $ws.Range("H2:H3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("H2:H4"))
}
foreach ($fc in $ws.Range("I2:I3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("I2:I4"))
}
foreach ($fc in $ws.Range("J2:J3").FormatConditions) {
    $fc.ModifyAppliesToRange($ws.Range("J2:J4"))
}

# Update the Dashboard summary count for "Opvolging / Status"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 3
